# Update "想去人数" (number of people interested) figures that were refreshed
# by the scraper run at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1872
$ws1.Range("F4").Value = 813
$ws1.Range("F5").Value = 679
$ws1.Range("F6").Value = 232

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1872
$ws4.Range("F5").Value = 813
$ws4.Range("F6").Value = 679
$ws4.Range("F7").Value = 232
